{"js": "// Remove the \"First sentence of US Declaration of Independence\" passage\n// (title paragraph, body paragraph, and the two blank paragraphs that\n// followed it), leaving \"First Paragraph of The Metamorphosis, Franz\n// Kafka\" as the new first passage in the document.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Identify the paragraphs that make up the Declaration of Independence\n// passage by matching on their text, rather than assuming a fixed\n// index, so the script is resilient to minor structural differences.\nconst targetTexts = [\n  \"First sentence of US Declaration of Independence\",\n  \"When in the Course of human events, it becomes necessary for one people to dissolve the political bands which have connected them with another, and to assume among the powers of the earth, the separate and equal station to which the Laws of Nature and of Nature's God entitle them, a decent respect to the opinions of mankind requires that they should declare the causes which impel them to the separation.\"\n];\n\nconst items = paragraphs.items;\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // The passage consists of the title paragraph, the quoted-text\n  // paragraph, and the two empty paragraphs that separate it from the\n  // next passage.\n  let endIndex = startIndex + 1;\n  if (endIndex < items.length && items[endIndex].text === targetTexts[1]) {\n    endIndex++;\n  }\n  while (endIndex < items.length && items[endIndex].text === \"\") {\n    endIndex++;\n  }\n\n  // Delete paragraphs from endIndex-1 down to startIndex (reverse order\n  // keeps earlier indices valid while deleting).\n  for (let i = endIndex - 1; i >= startIndex; i--) {\n    items[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"First sentence of US Declaration of Independence\" passage\n# (title paragraph, body paragraph, and the two blank paragraphs that\n# followed it), leaving \"First Paragraph of The Metamorphosis, Franz\n# Kafka\" as the new first passage in the document.\n\n$d = $word.ActiveDocument\n\n$titleText = \"First sentence of US Declaration of Independence\"\n$bodyText = \"When in the Course of human events, it becomes necessary for one people to dissolve the political bands which have connected them with another, and to assume among the powers of the earth, the separate and equal station to which the Laws of Nature and of Nature's God entitle them, a decent respect to the opinions of mankind requires that they should declare the causes which impel them to the separation.\"\n\n$count = $d.Paragraphs.Count\n$startIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($t -eq $titleText) {\n        $startIdx = $i\n        break\n    }\n}\n\nif ($startIdx -ne -1) {\n    $endIdx = $startIdx\n\n    # Include the following paragraph if it is the quoted passage text.\n    if (($startIdx + 1) -le $count) {\n        $nextText = $d.Paragraphs.Item($startIdx + 1).Range.Text.TrimEnd()\n        if ($nextText -eq $bodyText) {\n            $endIdx = $startIdx + 1\n        }\n    }\n\n    # Swallow any blank paragraphs that separate this passage from the\n    # next one.\n    while (($endIdx + 1) -le $count -and $d.Paragraphs.Item($endIdx + 1).Range.Text.TrimEnd() -eq \"\") {\n        $endIdx = $endIdx + 1\n    }\n\n    $startRange = $d.Paragraphs.Item($startIdx).Range\n    $endRange = $d.Paragraphs.Item($endIdx).Range\n    $deleteRange = $d.Range($startRange.Start, $endRange.End)\n    $deleteRange.Delete()\n}\n"}
